$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status of rows 11-14 (Post Numbers 10-13) from "Edited" to "Posted"
$ws.Range("D11").Value = "Posted"
$ws.Range("D12").Value = "Posted"
$ws.Range("D13").Value = "Posted"
$ws.Range("D14").Value = "Posted"

# Add new row 15 for Post Number 14 - "The Bahea"
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "The Bahea"
$ws.Range("C15").Value = "TheBahea.jpg"
$ws.Range("D15").Value = "Edited"

# Update selection to reflect the new active cell after edits
$ws.Range("G22").Select()
